$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 187, pushing the existing rows 187-190 down to 188-191
$ws.Rows.Item(187).Insert()

# Populate the newly inserted row 187 with the new weekly record.
# Non-numeric / repeated columns mirror the surrounding rows for this market/product.
$ws.Cells.Item(187, 1).Value = 7
$ws.Cells.Item(187, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(187, 3).Value = "Ñuble"
$ws.Cells.Item(187, 4).Value = 44595
$ws.Cells.Item(187, 5).Value = 16
$ws.Cells.Item(187, 6).Value = 100112003
$ws.Cells.Item(187, 7).Value = "Ajo"
$ws.Cells.Item(187, 8).Value = "Chino"
$ws.Cells.Item(187, 9).Value = "Primera"
$ws.Cells.Item(187, 10).Value = 100
$ws.Cells.Item(187, 11).Value = 19000
$ws.Cells.Item(187, 12).Value = 20000
$ws.Cells.Item(187, 13).Value = 19500
$ws.Cells.Item(187, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(187, 15).Value = "China"
$ws.Cells.Item(187, 16).Value = 1950
$ws.Cells.Item(187, 17).Value = 10
$ws.Cells.Item(187, 18).Value = "Hortaliza"
